$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NIT values for rows 2-11 (row 2 is updated, rows 3-11 are new)
$values = @(900687420, 900291664, 830018673, 900387627, 900123743, 901001861, 900003261, 900377532, 900190116, 900155802)

# Update A2 value first
$ws.Range("A2").Value = $values[0]

# Remove the fill ("No Fill") from A2's style so it matches the cleaned-up
# formatting (fillId -> 0, applyFill dropped) used by the migration sheet.
$ws.Range("A2").Interior.ColorIndex = -4142
$ws.Range("A2").Interior.Pattern = -4142

# Populate the new rows (3-11), copying A2's formatting (border/font/alignment)
# onto each new cell so they share the very same cell style as A2.
for ($i = 1; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("A2").Copy()
    $target = $ws.Cells.Item($row, 1)
    $target.PasteSpecial(-4122)
    $target.Value = $values[$i]
    $ws.Rows($row).RowHeight = 16
}

$excel.CutCopyMode = $false

# Match the selection left behind by the edit
$null = $ws.Range("A2:A11").Select()
